# Apply Kujata_Profits value updates (scheduled runner sync)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 15000
$ws.Range("J3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("N3").Value = -15228

$ws.Range("H33").Value = 450.75
$ws.Range("I33").Value = 419.42856
$ws.Range("J33").Value = 544.7143
$ws.Range("K33").Value = 419.42856
$ws.Range("L33").Value = 544.7143
$ws.Range("M33").Value = -190.42856
$ws.Range("N33").Value = -1002.7143

$ws.Range("H41").Value = 2817.75
$ws.Range("I41").Value = 2292.4
$ws.Range("J41").Value = 3693.3333
$ws.Range("K41").Value = 2292.4
$ws.Range("L41").Value = 3693.3333
$ws.Range("M41").Value = -1852.4
$ws.Range("N41").Value = -4573.3333

$ws.Range("H82").Value = 393.33334
$ws.Range("I82").Value = 393.33334
$ws.Range("K82").Value = 1180.00002
$ws.Range("M82").Value = -774.0000199999999

$ws.Range("H85").Value = 393.33334
$ws.Range("I85").Value = 393.33334
$ws.Range("K85").Value = 1180.00002
$ws.Range("M85").Value = 223.9999800000001

$ws.Range("H102").Value = 15000
$ws.Range("J102").Value = 15000
$ws.Range("L102").Value = 15000
$ws.Range("N102").Value = -21490

$ws.Range("H129").Value = 818.0303
$ws.Range("J129").Value = 927
$ws.Range("L129").Value = 2781
$ws.Range("N129").Value = -12781

$ws.Range("H138").Value = 459332.25
$ws.Range("J138").Value = 513159.75
$ws.Range("L138").Value = 1539479.25
$ws.Range("N138").Value = -1549759.25

$ws.Range("H141").Value = 4640.25
$ws.Range("I141").Value = 6330.5
$ws.Range("K141").Value = 18991.5
$ws.Range("M141").Value = -13811.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5197.143
$ws.Range("I31").Value = 5197.143
$ws.Range("K31").Value = 5197.143
$ws.Range("M31").Value = -4903.143

$ws.Range("H61").Value = 55556652
$ws.Range("I61").Value = 71429420
$ws.Range("J61").Value = 1944.75
$ws.Range("K61").Value = 71429420
$ws.Range("L61").Value = 1944.75
$ws.Range("M61").Value = -71429208
$ws.Range("N61").Value = -2368.75

$ws.Range("H74").Value = 2753.8
$ws.Range("I74").Value = 755
$ws.Range("K74").Value = 755
$ws.Range("M74").Value = 119

$ws.Range("H77").Value = 2753.8
$ws.Range("I77").Value = 755
$ws.Range("K77").Value = 3775
$ws.Range("M77").Value = 593

$ws.Range("H136").Value = 55556652
$ws.Range("I136").Value = 71429420
$ws.Range("J136").Value = 1944.75
$ws.Range("K136").Value = 214288260
$ws.Range("L136").Value = 5834.25
$ws.Range("M136").Value = -214285710
$ws.Range("N136").Value = -10934.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 59666.668
$ws.Range("J87").Value = 59666.668
$ws.Range("L87").Value = 59666.668
$ws.Range("N87").Value = -62162.668

$ws.Range("H90").Value = 59666.668
$ws.Range("J90").Value = 59666.668
$ws.Range("L90").Value = 179000.004
$ws.Range("N90").Value = -191480.004

$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1375.2245
$ws.Range("I31").Value = 1445.65
$ws.Range("J31").Value = 1326.6552
$ws.Range("K31").Value = 1445.65
$ws.Range("L31").Value = 1326.6552
$ws.Range("M31").Value = -1150.65
$ws.Range("N31").Value = -1916.6552

$ws.Range("H34").Value = 1375.2245
$ws.Range("I34").Value = 1445.65
$ws.Range("J34").Value = 1326.6552
$ws.Range("K34").Value = 1445.65
$ws.Range("L34").Value = 1326.6552
$ws.Range("M34").Value = -1243.65
$ws.Range("N34").Value = -1730.6552

$ws.Range("H134").Value = 16130768
$ws.Range("I134").Value = 1726.6296
$ws.Range("K134").Value = 5179.8888
$ws.Range("M134").Value = -2644.8888

$ws.Range("H141").Value = 1014264
$ws.Range("J141").Value = 1014264
$ws.Range("L141").Value = 1014264
$ws.Range("N141").Value = -1024624

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 700.4
$ws.Range("J29").Value = 850.5
$ws.Range("L29").Value = 2551.5
$ws.Range("N29").Value = -3105.5

$ws.Range("H64").Value = 4212.8
$ws.Range("I64").Value = 1264
$ws.Range("J64").Value = 4950
$ws.Range("K64").Value = 3792
$ws.Range("L64").Value = 14850
$ws.Range("M64").Value = -3522
$ws.Range("N64").Value = -15390

$ws.Range("H67").Value = 4212.8
$ws.Range("I67").Value = 1264
$ws.Range("J67").Value = 4950
$ws.Range("K67").Value = 3792
$ws.Range("L67").Value = 14850
$ws.Range("M67").Value = -2856
$ws.Range("N67").Value = -16722

$ws.Range("H70").Value = 6109
$ws.Range("J70").Value = 6320
$ws.Range("L70").Value = 18960
$ws.Range("N70").Value = -19590

$ws.Range("H73").Value = 6109
$ws.Range("J73").Value = 6320
$ws.Range("L73").Value = 18960
$ws.Range("N73").Value = -21144

$ws.Range("H75").Value = 2337.6667
$ws.Range("I75").Value = 1006.5
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 3019.5
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -2021.5
$ws.Range("N75").Value = -16996

$ws.Range("H78").Value = 2337.6667
$ws.Range("I78").Value = 1006.5
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 9058.5
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -4066.5
$ws.Range("N78").Value = -54984

$ws.Range("H117").Value = 852.93335
$ws.Range("I117").Value = 444
$ws.Range("J117").Value = 1210.75
$ws.Range("K117").Value = 1332
$ws.Range("L117").Value = 3632.25
$ws.Range("M117").Value = 2110
$ws.Range("N117").Value = -10516.25

$ws.Range("H129").Value = 20834456
$ws.Range("I129").Value = 41667320
$ws.Range("J129").Value = 6945878
$ws.Range("K129").Value = 125001960
$ws.Range("L129").Value = 20837634
$ws.Range("M129").Value = -124996960
$ws.Range("N129").Value = -20847634

$ws.Range("H131").Value = 32259878
$ws.Range("I131").Value = 125000350
$ws.Range("J131").Value = 2322.348
$ws.Range("K131").Value = 375001050
$ws.Range("L131").Value = 6967.044
$ws.Range("M131").Value = -374996010
$ws.Range("N131").Value = -17047.044

$ws.Range("H139").Value = 1594.6923
$ws.Range("I139").Value = 1536.72
$ws.Range("K139").Value = 4610.16
$ws.Range("M139").Value = 529.8400000000001

$ws.Range("H140").Value = 22264.27
$ws.Range("I140").Value = 52449.35
$ws.Range("K140").Value = 157348.05
$ws.Range("M140").Value = -152168.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3848.1428
$ws.Range("I80").Value = 3417.0833
$ws.Range("J80").Value = 4422.8887
$ws.Range("K80").Value = 3417.0833
$ws.Range("L80").Value = 4422.8887
$ws.Range("M80").Value = -2419.0833
$ws.Range("N80").Value = -6418.8887

$ws.Range("H83").Value = 3848.1428
$ws.Range("I83").Value = 3417.0833
$ws.Range("J83").Value = 4422.8887
$ws.Range("K83").Value = 17085.4165
$ws.Range("L83").Value = 22114.4435
$ws.Range("M83").Value = -12093.4165
$ws.Range("N83").Value = -32098.4435

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8600
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 8600
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 8600
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -8976

$ws.Range("H122").Value = 25002020
$ws.Range("I122").Value = 35716144
$ws.Range("K122").Value = 107148432
$ws.Range("M122").Value = -107145982

$ws.Range("H136").Value = 1770.1428
$ws.Range("I136").Value = 1278.3
$ws.Range("J136").Value = 2999.75
$ws.Range("K136").Value = 3834.9
$ws.Range("L136").Value = 8999.25
$ws.Range("M136").Value = -1284.9
$ws.Range("N136").Value = -14099.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2026.1936
$ws.Range("I132").Value = 1837.5652
$ws.Range("J132").Value = 2568.5
$ws.Range("K132").Value = 5512.6956
$ws.Range("L132").Value = 7705.5
$ws.Range("M132").Value = -2982.6956
$ws.Range("N132").Value = -12765.5

$ws.Range("H136").Value = 1267.96
$ws.Range("I136").Value = 1263.421
$ws.Range("J136").Value = 1282.3334
$ws.Range("K136").Value = 3790.263
$ws.Range("L136").Value = 3847.0002
$ws.Range("M136").Value = -1240.263
$ws.Range("N136").Value = -8947.0002
